$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in cell A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 08:46"

# Row 42 - Singapur
$ws.Cells.Item(42, 2).Value = 45613
$ws.Cells.Item(42, 5).Value = 3807

# Row 48 - Afganistan
$ws.Cells.Item(48, 2).Value = 34351
$ws.Cells.Item(48, 3).Value = 157
$ws.Cells.Item(48, 4).Value = 21105
$ws.Cells.Item(48, 5).Value = 12271
$ws.Cells.Item(48, 7).Value = 4
$ws.Cells.Item(48, 8).Value = 975

# Row 76 - El Salvador
$ws.Cells.Item(76, 4).Value = 5454
$ws.Cells.Item(76, 5).Value = 3434
$ws.Cells.Item(76, 7).Value = 5
$ws.Cells.Item(76, 8).Value = 254

# Row 98 - Hungria
$ws.Cells.Item(98, 2).Value = 4229
$ws.Cells.Item(98, 3).Value = 6
$ws.Cells.Item(98, 4).Value = 2974
$ws.Cells.Item(98, 5).Value = 660
$ws.Cells.Item(98, 7).Value = 2
$ws.Cells.Item(98, 8).Value = 595

# Row 143 - Georgia
$ws.Cells.Item(143, 2).Value = 981
$ws.Cells.Item(143, 3).Value = 8
$ws.Cells.Item(143, 4).Value = 851
$ws.Cells.Item(143, 5).Value = 115

# Row 178 - Camboya
$ws.Cells.Item(178, 4).Value = 133
$ws.Cells.Item(178, 5).Value = 8
